$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds pre-formatted price strings (e.g. "43.531.10", "5.70") that
# look numeric to Excel and would otherwise be auto-coerced to a Double,
# dropping the literal text (trailing zeros, European grouping dots, etc).
# Force the target range to Text before writing, then restore the default
# "Normal" style afterwards so the saved cell style matches the original
# (no explicit style index).
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '43.531.10'
$ws.Range("E2").Value = '  +2.62%  '
$ws.Range("D3").Value = '2.419.85'
$ws.Range("E3").Value = '  +8.84%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '325.67'
$ws.Range("E5").Value = '  +13.16%  '
$ws.Range("D6").Value = '104.91'
$ws.Range("E6").Value = '  -4.59%  '
$ws.Range("D7").Value = '0.646'
$ws.Range("E7").Value = '  +4.03%  '
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("E9").Value = '  +10.65%  '
$ws.Range("D10").Value = '42.14'
$ws.Range("E10").Value = '  -2.55%  '
$ws.Range("D11").Value = '0.0951'
$ws.Range("E11").Value = '  +4.70%  '
$ws.Range("D12").Value = '8.68'
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("E13").Value = '  +3.33%  '
$ws.Range("D14").Value = '17.28'
$ws.Range("E14").Value = '  +16.45%  '
$ws.Range("E15").Value = '  +2.77%  '
$ws.Range("D16").Value = '2.777.67'
$ws.Range("E16").Value = '  +8.58%  '
$ws.Range("D17").Value = '2.492.86'
$ws.Range("E17").Value = '  +11.49%  '
$ws.Range("D18").Value = '43.543.00'
$ws.Range("E18").Value = '  +2.76%  '
$ws.Range("E19").Value = '  +6.07%  '
$ws.Range("D20").Value = '7.49'
$ws.Range("E20").Value = '  +5.10%  '
$ws.Range("D21").Value = '75.62'
$ws.Range("E21").Value = '  +3.63%  '
$ws.Range("D22").Value = '3.52'
$ws.Range("E22").Value = '  +4.24%  '
$ws.Range("D23").Value = '260.97'
$ws.Range("E23").Value = '  +13.56%  '
$ws.Range("E24").Value = '  +3.37%  '
$ws.Range("D25").Value = '9.58'
$ws.Range("E25").Value = '  +7.41%  '
$ws.Range("D26").Value = '12.01'
$ws.Range("E26").Value = '  +5.62%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  +0.32%  '
$ws.Range("D29").Value = '22.95'
$ws.Range("E29").Value = '  +10.63%  '
$ws.Range("D30").Value = '179.84'
$ws.Range("E30").Value = '  +3.98%  '
$ws.Range("E31").Value = '  +2.04%  '
$ws.Range("D32").Value = '38.32'
$ws.Range("E32").Value = '  +3.73%  '
$ws.Range("E33").Value = '  +2.42%  '
$ws.Range("D34").Value = '0.0937'
$ws.Range("E34").Value = '  +7.41%  '
$ws.Range("D35").Value = '5.95'
$ws.Range("E35").Value = '  +6.91%  '
$ws.Range("E36").Value = '  +6.13%  '
$ws.Range("D37").Value = '4.91'
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("D38").Value = '0.0372'
$ws.Range("E38").Value = '  +1.47%  '
$ws.Range("D39").Value = '4.01'
$ws.Range("E39").Value = '  -3.42%  '
$ws.Range("D40").Value = '0.107'
$ws.Range("E40").Value = '  +3.19%  '
$ws.Range("D41").Value = '2.92'
$ws.Range("E41").Value = '  +23.40%  '
$ws.Range("D42").Value = '1.63'
$ws.Range("E42").Value = '  +26.50%  '
$ws.Range("D43").Value = '0.235'
$ws.Range("E43").Value = '  +2.62%  '
$ws.Range("D44").Value = '125.54'
$ws.Range("E44").Value = '  +24.58%  '
$ws.Range("D45").Value = '69.85'
$ws.Range("E45").Value = '  -5.04%  '
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").Value = '12.65'
$ws.Range("E47").Value = '  +3.33%  '
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").Value = '5.70'
$ws.Range("E48").Value = '  +7.32%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '9.52'
$ws.Range("E49").Value = '  +13.39%  '
$ws.Range("E50").Value = '  +4.19%  '
$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").Value = '75.35'
$ws.Range("E51").Value = '  +11.90%  '

$priceRange.Style = "Normal"
